# Weekly update: insert this week's new price record at the top of the
# data table (row 2), pushing all existing records down by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new blank row right above the current first data row.
$ws.Rows("2:2").Insert()

# The inserted row inherits the header row's formatting; strip that back
# to the plain "data row" look used by every other record.
$ws.Rows("2:2").ClearFormats()

# Column D carries the date number format - copy it from the row below
# (an existing, correctly-formatted data row) rather than hard-coding a
# style index.
$ws.Range("D2").NumberFormat = $ws.Range("D3").NumberFormat

# Populate the new record's values.
$ws.Range("A2").Value = 8
$ws.Range("B2").Value = "Terminal La Palmera de La Serena"
$ws.Range("C2").Value = "Coquimbo"
$ws.Range("D2").Value = 44599
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = 100114007
$ws.Range("G2").Value = "Jengibre"
$ws.Range("H2").Value = "Sin especificar"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 400
$ws.Range("K2").Value = 15000
$ws.Range("L2").Value = 16000
$ws.Range("M2").Value = 15500
$ws.Range("N2").Value = '$/caja 13 kilos'
$ws.Range("O2").Value = "Perú"
$ws.Range("P2").Value = 1192
$ws.Range("Q2").Value = 13
$ws.Range("R2").Value = "Hortaliza"
